# Applies the "Add files via upload" edit: appends a line break to the
# final "If no PCB exists..." list item, then appends a new block of
# numbered-list content (test-process loader notes + alarm command notes)
# at the end of the document body, matching the same ListParagraph /
# numId=3 bullet list used throughout the manual.

$d = $word.ActiveDocument

function Add-ListItem($level, $parts) {
    # Inserts a brand-new paragraph at the end of the document, set to the
    # requested outline level (0 = top-level bullet, 1 = sub-bullet) of the
    # existing numId=3 list, then types each element of $parts into it.
    # A part that is the literal token "@BR@" inserts a manual line break
    # instead of text.
    $tail = $d.Paragraphs.Last
    $tail.Range.InsertParagraphAfter() | Out-Null
    $p = $d.Paragraphs.Last
    $p.Range.ListFormat.ListLevelNumber = $level + 1

    foreach ($part in $parts) {
        $r = $p.Range
        $r.Collapse(0)  # wdCollapseEnd
        if ($part -eq "@BR@") {
            $r.InsertBreak(6)  # wdLineBreak
        } else {
            $r.InsertAfter($part)
        }
    }
    return $p
}

# --- finish off the existing last paragraph with a manual line break ---
$last = $d.Paragraphs.Last
$br = $last.Range
$br.Collapse(0)  # wdCollapseEnd
$br.InsertBreak(6)  # wdLineBreak

# --- "To load test processes: loadr3" block ---
Add-ListItem 0 @("To load test processes: loadr3") | Out-Null
Add-ListItem 1 @("Loads a total of five test processes in a suspended ready state.") | Out-Null
Add-ListItem 1 @("When a test process is unsuspended, it will print a message indicating its successful run in the CPU.") | Out-Null
Add-ListItem 1 @("A test process will perform its task more than once based on the process" + [char]8217 + " number.") | Out-Null
Add-ListItem 1 @("Test processes can be completely removed from the system before their termination as long as the process is in a non-suspended state. Otherwise, the process cannot be removed", " prematurely.", "@BR@") | Out-Null

# --- "To create an alarm: alarm [message] [time]" block ---
Add-ListItem 0 @("To ", " ", "create an alarm: alarm [message] [time]") | Out-Null
Add-ListItem 1 @("Creates an alarm that will print a given message to the screen when a specified time is reached.") | Out-Null
Add-ListItem 1 @("[message] is the message that the alarm will print to the screen.") | Out-Null
Add-ListItem 1 @("[time] is the amount of time in seconds that must at least pass from the time of creating the alarm before the message can be printed.", " This value should a positive integer.") | Out-Null
Add-ListItem 1 @("Depending on the CPU load and the number of processes preceding the alarm-checking process, it is possible for the message to be printed later than the specified time.") | Out-Null
Add-ListItem 1 @("It is possible to create multiple alarms with different or similar messages which can be printed together given that the timer for each alarm had been reached.") | Out-Null

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
